$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 19: new time-tracking entry (9.12) ---
$ws.Range("A19").Value = 9.12
$ws.Range("B19").Value = 0.79166666666666663
$ws.Range("B19").NumberFormat = "h:mm"
$ws.Range("C19").Value = 0.875
$ws.Range("C19").NumberFormat = "h:mm"
$ws.Range("E19").Value = "2hr"
$ws.Range("F19").Value = "adding global navigation as famous website's logo and trim little about number chart function"

# --- Row 20: new time-tracking entry (10.12) ---
$ws.Range("A20").Value = 10.12
$ws.Range("B20").Value = 0.91666666666666663
$ws.Range("B20").NumberFormat = "h:mm"
$ws.Range("C20").Value = 1
$ws.Range("C20").NumberFormat = "[h]:mm:ss"
$ws.Range("E20").Value = "2hr"
$ws.Range("F20").Value = "finding out error in random values function so try to fix it"

# --- Update the active selection to reflect where the user left off ---
$ws.Range("B21").Select()
